$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F (match the bold/border/centered style used
# by the other header cells B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160

# Timestamps for rows 2..207 (data rows), matching the commit's time_taken values
$times = @(
"2021-10-05 13:38:34.270574",
"2021-10-05 13:38:34.270585",
"2021-10-05 13:38:34.270589",
"2021-10-05 13:38:34.270591",
"2021-10-05 13:38:34.270594",
"2021-10-05 13:38:34.270597",
"2021-10-05 13:38:34.270599",
"2021-10-05 13:38:34.270602",
"2021-10-05 13:38:34.270604",
"2021-10-05 13:38:34.270607",
"2021-10-05 13:38:34.270609",
"2021-10-05 13:38:34.270612",
"2021-10-05 13:38:34.270614",
"2021-10-05 13:38:34.270617",
"2021-10-05 13:38:34.270619",
"2021-10-05 13:38:34.270621",
"2021-10-05 13:38:34.270624",
"2021-10-05 13:38:34.270627",
"2021-10-05 13:38:34.270629",
"2021-10-05 13:38:34.270632",
"2021-10-05 13:38:34.270634",
"2021-10-05 13:38:34.270636",
"2021-10-05 13:38:34.270639",
"2021-10-05 13:38:34.270641",
"2021-10-05 13:38:34.270648",
"2021-10-05 13:38:34.270651",
"2021-10-05 13:38:34.270653",
"2021-10-05 13:38:34.270656",
"2021-10-05 13:38:34.270658",
"2021-10-05 13:38:34.270661",
"2021-10-05 13:38:34.270663",
"2021-10-05 13:38:34.270665",
"2021-10-05 13:38:34.270668",
"2021-10-05 13:38:34.270671",
"2021-10-05 13:38:34.270673",
"2021-10-05 13:38:34.270676",
"2021-10-05 13:38:34.270678",
"2021-10-05 13:38:34.270680",
"2021-10-05 13:38:34.270683",
"2021-10-05 13:38:34.270685",
"2021-10-05 13:38:34.270688",
"2021-10-05 13:38:34.270690",
"2021-10-05 13:38:34.270693",
"2021-10-05 13:38:34.270695",
"2021-10-05 13:38:34.270698",
"2021-10-05 13:38:34.270700",
"2021-10-05 13:38:34.270703",
"2021-10-05 13:38:34.270705",
"2021-10-05 13:38:34.270708",
"2021-10-05 13:38:34.270710",
"2021-10-05 13:38:34.270713",
"2021-10-05 13:38:34.270715",
"2021-10-05 13:38:34.270718",
"2021-10-05 13:38:34.270721",
"2021-10-05 13:38:34.270723",
"2021-10-05 13:38:34.270725",
"2021-10-05 13:38:34.270728",
"2021-10-05 13:38:34.270730",
"2021-10-05 13:38:34.270733",
"2021-10-05 13:38:34.270735",
"2021-10-05 13:38:34.270738",
"2021-10-05 13:38:34.270740",
"2021-10-05 13:38:34.270743",
"2021-10-05 13:38:34.270745",
"2021-10-05 13:38:34.270749",
"2021-10-05 13:38:34.270751",
"2021-10-05 13:38:34.270754",
"2021-10-05 13:38:34.270756",
"2021-10-05 13:38:34.270758",
"2021-10-05 13:38:34.270761",
"2021-10-05 13:38:34.270763",
"2021-10-05 13:38:34.270766",
"2021-10-05 13:38:34.270768",
"2021-10-05 13:38:34.270771",
"2021-10-05 13:38:34.270773",
"2021-10-05 13:38:34.270775",
"2021-10-05 13:38:34.270779",
"2021-10-05 13:38:34.270782",
"2021-10-05 13:38:34.270785",
"2021-10-05 13:38:34.270787",
"2021-10-05 13:38:34.270790",
"2021-10-05 13:38:34.270792",
"2021-10-05 13:38:34.270795",
"2021-10-05 13:38:34.270797",
"2021-10-05 13:38:34.270799",
"2021-10-05 13:38:34.270802",
"2021-10-05 13:38:34.270804",
"2021-10-05 13:38:34.270807",
"2021-10-05 13:38:34.270809",
"2021-10-05 13:38:34.270811",
"2021-10-05 13:38:34.270814",
"2021-10-05 13:38:34.270816",
"2021-10-05 13:38:34.270820",
"2021-10-05 13:38:34.270823",
"2021-10-05 13:38:34.270825",
"2021-10-05 13:38:34.270827",
"2021-10-05 13:38:34.270830",
"2021-10-05 13:38:34.270832",
"2021-10-05 13:38:34.270835",
"2021-10-05 13:38:34.270837",
"2021-10-05 13:38:34.270840",
"2021-10-05 13:38:34.270842",
"2021-10-05 13:38:34.270844",
"2021-10-05 13:38:34.270847",
"2021-10-05 13:38:34.270849",
"2021-10-05 13:38:34.270852",
"2021-10-05 13:38:34.270854",
"2021-10-05 13:38:34.270856",
"2021-10-05 13:38:34.270861",
"2021-10-05 13:38:34.270864",
"2021-10-05 13:38:34.270866",
"2021-10-05 13:38:34.270869",
"2021-10-05 13:38:34.270871",
"2021-10-05 13:38:34.270873",
"2021-10-05 13:38:34.270876",
"2021-10-05 13:38:34.270878",
"2021-10-05 13:38:34.270881",
"2021-10-05 13:38:34.270883",
"2021-10-05 13:38:34.270886",
"2021-10-05 13:38:34.270888",
"2021-10-05 13:38:34.270890",
"2021-10-05 13:38:34.270893",
"2021-10-05 13:38:34.270895",
"2021-10-05 13:38:34.270900",
"2021-10-05 13:38:34.270903",
"2021-10-05 13:38:34.270905",
"2021-10-05 13:38:34.270908",
"2021-10-05 13:38:34.270910",
"2021-10-05 13:38:34.270914",
"2021-10-05 13:38:34.270917",
"2021-10-05 13:38:34.270920",
"2021-10-05 13:38:34.270922",
"2021-10-05 13:38:34.270924",
"2021-10-05 13:38:34.270927",
"2021-10-05 13:38:34.270929",
"2021-10-05 13:38:34.270932",
"2021-10-05 13:38:34.270934",
"2021-10-05 13:38:34.270937",
"2021-10-05 13:38:34.270939",
"2021-10-05 13:38:34.270942",
"2021-10-05 13:38:34.270944",
"2021-10-05 13:38:34.270946",
"2021-10-05 13:38:34.270949",
"2021-10-05 13:38:34.270951",
"2021-10-05 13:38:34.270954",
"2021-10-05 13:38:34.270956",
"2021-10-05 13:38:34.270958",
"2021-10-05 13:38:34.270961",
"2021-10-05 13:38:34.270964",
"2021-10-05 13:38:34.270966",
"2021-10-05 13:38:34.270969",
"2021-10-05 13:38:34.270971",
"2021-10-05 13:38:34.270973",
"2021-10-05 13:38:34.270976",
"2021-10-05 13:38:34.270978",
"2021-10-05 13:38:34.270981",
"2021-10-05 13:38:34.270983",
"2021-10-05 13:38:34.270985",
"2021-10-05 13:38:34.270988",
"2021-10-05 13:38:34.270990",
"2021-10-05 13:38:34.270993",
"2021-10-05 13:38:34.270995",
"2021-10-05 13:38:34.270998",
"2021-10-05 13:38:34.271000",
"2021-10-05 13:38:34.271002",
"2021-10-05 13:38:34.271005",
"2021-10-05 13:38:34.271007",
"2021-10-05 13:38:34.271010",
"2021-10-05 13:38:34.271012",
"2021-10-05 13:38:34.271014",
"2021-10-05 13:38:34.271017",
"2021-10-05 13:38:34.271019",
"2021-10-05 13:38:34.271023",
"2021-10-05 13:38:34.271028",
"2021-10-05 13:38:34.271031",
"2021-10-05 13:38:34.271033",
"2021-10-05 13:38:34.271036",
"2021-10-05 13:38:34.271038",
"2021-10-05 13:38:34.271041",
"2021-10-05 13:38:34.271043",
"2021-10-05 13:38:34.271046",
"2021-10-05 13:38:34.271048",
"2021-10-05 13:38:34.271050",
"2021-10-05 13:38:34.271053",
"2021-10-05 13:38:34.271055",
"2021-10-05 13:38:34.271058",
"2021-10-05 13:38:34.271060",
"2021-10-05 13:38:34.271062",
"2021-10-05 13:38:34.271065",
"2021-10-05 13:38:34.271067",
"2021-10-05 13:38:34.271070",
"2021-10-05 13:38:34.271072",
"2021-10-05 13:38:34.271075",
"2021-10-05 13:38:34.271077",
"2021-10-05 13:38:34.271080",
"2021-10-05 13:38:34.271082",
"2021-10-05 13:38:34.271084",
"2021-10-05 13:38:34.271087",
"2021-10-05 13:38:34.271089",
"2021-10-05 13:38:34.271092",
"2021-10-05 13:38:34.271094",
"2021-10-05 13:38:34.271097",
"2021-10-05 13:38:34.271099",
"2021-10-05 13:38:34.271102",
"2021-10-05 13:38:34.271104",
"2021-10-05 13:38:34.271107"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
